# OSLC change set delivery.pptx - apply commit "Jra change set delivery (#620)"
# 1) Refresh cached date fields (2/5/24 -> 8/16/24) on the Handout Master,
#    Notes Master, and the 12 unused "bank" slide layouts (idx 32-43 on
#    SlideMaster1's CustomLayouts) that still carry the old cached date.
# 2) Split the "delivery is effected by ..." bullet on slide 7 into three
#    runs, correcting "effected" -> "affected" in the process.
# 3) Minor touch-up of the trailing empty placeholder on slide 7.

$p = $ppt.ActivePresentation

$oldDate = "2/5/24"
$newDate = "8/16/24"

# --- Handout Master date field -------------------------------------------------
$hm = $p.HandoutMaster
$hmDate = $hm.HeadersFooters.DateAndTime
$hmDate.Text = $newDate

# --- Notes Master date field -----------------------------------------------------
$nm = $p.NotesMaster
$nmDate = $nm.HeadersFooters.DateAndTime
$nmDate.Text = $newDate

# --- Unused bank slide layouts (CustomLayouts 32-43 on SlideMaster 1) ------------
$customLayouts = $p.SlideMaster.CustomLayouts
$layoutIndexes = 32,33,34,35,36,37,38,39,40,41,42,43
foreach ($idx in $layoutIndexes) {
    $cl = $customLayouts.Item($idx)
    $dateShape = $cl.Shapes.Item("Date Placeholder 4")
    if ($dateShape.TextFrame.TextRange.Text -eq $oldDate) {
        $dateShape.TextFrame.TextRange.Text = $newDate
    }
}

# --- Slide 7: fix bullet text and split run -----------------------------------
$s7 = $p.Slides.Item(7)

$bulletShape = $s7.Shapes.Item("Text Placeholder 2")
$tr = $bulletShape.TextFrame.TextRange
$para = $tr.Paragraphs(2)
$paraStart = $para.Start

$part1 = "delivery "
$part2 = "is affected "
$part3 = "by creating a delivery rather than creating a transient delivery session and then doing a GET+PUT on it"

# Rewrite just the "is effected " segment -> "is affected " (offsets are within
# the original, unmodified paragraph text), which the COM host naturally
# splits into three sibling runs.
$midRun = $tr.Characters($paraStart + $part1.Length, $part2.Length)
$midRun.Text = $part2

# Note: slide 7's trailing empty placeholder ("Text Placeholder 3") only
# gains a `dirty="0"` marker on its endParaRPr in the target deck - that
# attribute is an internal PowerPoint bookkeeping flag with no COM-exposed
# property, and any attempt to touch that empty text range via the object
# model (even a no-op re-assignment) forces the host to materialise an
# empty <a:r> run and drop the endParaRPr entirely, which is a bigger
# structural deviation than simply leaving the shape untouched. So it is
# intentionally left alone here.
